# Update "Name of Algo" results - correct the imputed KNN values for the
# terrestrial_mammals / combination_2_ABCDE / ABC / 10 / seed5 dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.519
$ws.Range("C3").Value = -12.516
$ws.Range("A12").Value = -21.629
$ws.Range("B14").Value = 5.767
$ws.Range("B26").Value = 6.193000000000001
$ws.Range("C30").Value = -12.941
$ws.Range("B31").Value = 6.205
$ws.Range("A32").Value = -21.406
$ws.Range("B35").Value = 8.095000000000001
$ws.Range("A36").Value = -21.12
$ws.Range("B37").Value = 8.260000000000002
$ws.Range("A38").Value = -20.363
$ws.Range("C44").Value = -12.289
$ws.Range("B45").Value = 5.569000000000001
$ws.Range("A46").Value = -21.556
$ws.Range("A54").Value = -21.876
$ws.Range("A55").Value = -21.921
$ws.Range("B57").Value = 6.090000000000001
$ws.Range("C58").Value = -12.609
$ws.Range("A67").Value = -21.588
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.55
$ws.Range("C84").Value = -13.659
$ws.Range("C89").Value = -11.275
$ws.Range("A91").Value = -21.747
$ws.Range("C91").Value = -11.21
$ws.Range("C92").Value = -11.539
$ws.Range("A99").Value = -20.559
$ws.Range("B100").Value = 6.187
$ws.Range("B102").Value = 7.334999999999999
$ws.Range("C102").Value = -12.808
